$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 648 ("「笑ってくれ、絶望を癒すのは君の笑みだから」" entry),
# shifting all subsequent rows up by one.
$ws.Rows.Item(648).Delete()
